$wb = $excel.ActiveWorkbook

# Insert a new "Documentation" worksheet before the first (currently active)
# sheet, so it becomes the new first/active tab.
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "Documentation"

# Populate the documentation text describing the workbook and its sheets.
$ws.Range("A1").Value = "Workbook: Mapping and transformation of new process in Agriculture, Forestry and Fishing Future Technologies across regions  for Tui Scenario "
$ws.Range("A2").Value = "AVA: Process availability across regions"
$ws.Range("A3").Value = "AF_Trans: Time-slice-specific availability factors for existing energy processes"
$ws.Range("A4").Value = "FILL Table: Model data for transformation operations"

[void]$ws.Range("F17").Select()
